$d = $word.ActiveDocument

# --- Header / footer -------------------------------------------------
# Attach a default header and a default footer to the document's only
# section. Writing (empty) content into the header/footer Range is what
# actually mints the header1.xml / footer1.xml parts and wires up the
# <w:headerReference w:type="default".../> / <w:footerReference .../>
# entries on the section's sectPr.
$section = $d.Sections.Item(1)

$header = $section.Headers.Item(1)   # wdHeaderFooterPrimary
$header.Range.InsertAfter("")

$footer = $section.Footers.Item(1)   # wdHeaderFooterPrimary
$footer.Range.InsertAfter("")

# --- List / numbering definition -------------------------------------
# Mint a numbering part (abstractNum + num) describing the look of
# bulleted/numbered lists, same as applying the default bullet list
# style from the list gallery. We apply it to a scratch paragraph (so
# the numbering definition gets created) and then remove that scratch
# paragraph again so the body's visible content is left untouched.
$scratch = $d.Content
$scratch.Collapse(0)
$scratch.InsertParagraphAfter()
$listParagraph = $d.Paragraphs.Last
$listParagraph.Range.ListFormat.ApplyBulletDefault()
$listParagraph.Range.Delete()
